# The workbook's single sheet ("Analysis") gains one new column, inserted
# immediately before column E. Everything that used to live in columns
# E.. shifts one column to the right (E->F, F->G, ... T->U), and the new
# column E is populated with a label for each of the four header rows and
# with "pair-a" for every data row (this raw-data export apparently always
# recorded a single object pair, so the new "object.pair" column is
# constant across all 16 trials).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before column E; this shifts existing columns
# E:T right to F:U, carrying their values/types along untouched.
$ws.Columns("E:E").Insert()

# Row 1/2 repeat the "Independent Variable" banner across columns C:G
# after the insert, so the new E cell just repeats that label too.
$ws.Range("E1").Value = "Independent Variable"
$ws.Range("E2").Value = "Independent Variable"

# Row 3 is the real header row: the new column is "object.pair".
$ws.Range("E3").Value = "object.pair"

# Row 4 (units row) has no unit for this column - it must stay blank,
# but still exist as an actual (empty-string) cell rather than being
# absent, matching the rest of the blank unit cells in that row. Writing
# a bare "" does not materialize a cell, so force it via a leading
# quote-prefix entry and then strip the resulting quote-prefix style
# back to Normal.
$ws.Range("E4").Value = "'"
$ws.Range("E4").Style = "Normal"

# Data rows 5-20: every trial belongs to the same object pair.
for ($r = 5; $r -le 20; $r++) {
    $ws.Cells.Item($r, 5).Value = "pair-a"
}
